# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" across all
# sheets, and refresh the associated handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status column (Overview B2/C2, zh-cn C2, de-de C2) -> Ready for handoff
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest Handoff Datetime for zh-cn
$wsZhCn.Range("E2").Value = "2016-03-24 20:45:17"

# Latest Handoff Datetime for de-de, mirrored on the Overview's
# "Latest Handoff Date" column (both shared the same timestamp string).
$wsDeDe.Range("E2").Value = "2016-03-24 20:45:21"
$wsOverview.Range("D2").Value = "2016-03-24 20:45:21"
